# Insert a new price-record row at row 322 (pushing the existing rows
# 322-397 down to 323-398) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(322).Insert()

$ws.Cells.Item(322, 1).Value  = 3
$ws.Cells.Item(322, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(322, 3).Value  = "Coquimbo"
$ws.Cells.Item(322, 4).Value2 = 44782
$ws.Cells.Item(322, 5).Value  = 5
$ws.Cells.Item(322, 6).Value  = 100114013
$ws.Cells.Item(322, 7).Value  = "Zanahoria"
$ws.Cells.Item(322, 8).Value  = "Sin especificar"
$ws.Cells.Item(322, 9).Value  = "Primera"
$ws.Cells.Item(322, 10).Value = 200
$ws.Cells.Item(322, 11).Value = 11000
$ws.Cells.Item(322, 12).Value = 12000
$ws.Cells.Item(322, 13).Value = 11450
$ws.Cells.Item(322, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(322, 16).Value = 572
$ws.Cells.Item(322, 17).Value = 20
$ws.Cells.Item(322, 18).Value = "Hortaliza"
